$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("common_forms")

# The SAE-related rows 13-16 on the "common_forms" sheet need to be
# reordered: row 13 ("SAE Awareness date") moves down to row 16, and the
# three rows below it (SAE Start date, SAE End date, SAE Category) each
# shift up by one row. Only columns A (var) and C (item_name) differ
# between these rows; D and E are identical for all of them.
# Use a temporary holding row (well below the used range) to stash row 13
# before overwriting it, preserving cell formatting (e.g. the shaded
# C14 cell) via Copy.

$holdRow = 100

$ws.Range("A13").Copy($ws.Range("A" + $holdRow))
$ws.Range("C13").Copy($ws.Range("C" + $holdRow))

$ws.Range("A14").Copy($ws.Range("A13"))
$ws.Range("C14").Copy($ws.Range("C13"))

$ws.Range("A15").Copy($ws.Range("A14"))
$ws.Range("C15").Copy($ws.Range("C14"))

$ws.Range("A16").Copy($ws.Range("A15"))
$ws.Range("C16").Copy($ws.Range("C15"))

$ws.Range("A" + $holdRow).Copy($ws.Range("A16"))
$ws.Range("C" + $holdRow).Copy($ws.Range("C16"))

$ws.Range("A" + $holdRow).Clear()
$ws.Range("C" + $holdRow).Clear()

# Make "common_forms" the active/selected sheet (was "column_names").
$ws.Activate()
